# Update average_county_temperature (column K) with refreshed NOAA data,
# and recompute the dependent ASHP COP columns (R, S) for the
# electrified_utilities rows whose temperature changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ID facilities (rows 2, 9, 10, 11, 14) ---
$ws.Range("K2").Value  = 1.925925925925943
$ws.Range("K9").Value  = 1.925925925925943
$ws.Range("K10").Value = 1.925925925925943
$ws.Range("K11").Value = 1.925925925925943
$ws.Range("K14").Value = 1.925925925925943

# --- ME facility (row 12) ---
$ws.Range("K12").Value = -1.226851851851833

# --- FL facilities (rows 18, 19, 20) ---
$ws.Range("K18").Value = 13.17361111111111
$ws.Range("K19").Value = 13.17361111111111
$ws.Range("K20").Value = 13.17361111111111

# Dependent worst/best ASHP COP recalculations for the rows that carry an
# electrified_utilities option (row 10 -> ID, row 19 -> FL)
$ws.Range("R10").Value = 1.204711500590784
$ws.Range("S10").Value = 1.258324667221298

$ws.Range("R19").Value = 1.281341554412284
$ws.Range("S19").Value = 1.343162802314449
